# "keyboard finished project cleanup"
# Adds new Typography rows (Display, Keyboard, Mode, button_labels) and the
# corresponding Translation rows (keyboard text ids + a couple of other
# clean-up rows) that go with them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell and strip the style index the COM
# layer likes to stamp onto freshly-materialised cells, so the cell ends
# up using the column's implicit/default style (no explicit s="...").
# Passing $null still materialises the (then totally empty) cell, which
# is what the target file expects for the untouched columns in each row.
# ---------------------------------------------------------------------
function Set-PlainCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Same as Set-PlainCell, but forces the cell to stay text even when the
# string looks like a number (e.g. "123"), which Excel would otherwise
# silently coerce into a numeric cell.
function Set-PlainTextCell($ws, $row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wildcardChars = '!”"#*%&()''$+-@_, .:;?/~±×÷•º`´{}©£€^®¥_=[]¡¢|\¿><'
$placeholder = '<placeholder>'

# ----------------------- Typography sheet -----------------------------
$tg = $wb.Worksheets.Item("Typography")

# Row 9 - Display
Set-PlainCell $tg 9 2 "Display"
Set-PlainCell $tg 9 3 "RobotoCondensed-Regular.ttf"
Set-PlainCell $tg 9 4 28
Set-PlainCell $tg 9 5 4
Set-PlainCell $tg 9 6 "?"
Set-PlainCell $tg 9 7 $wildcardChars
Set-PlainCell $tg 9 8 $null
Set-PlainCell $tg 9 9 "a-z,A-Z,0-9"
Set-PlainCell $tg 9 10 $null

# Row 10 - Keyboard
Set-PlainCell $tg 10 2 "Keyboard"
Set-PlainCell $tg 10 3 "RobotoCondensed-Regular.ttf"
Set-PlainCell $tg 10 4 20
Set-PlainCell $tg 10 5 4
Set-PlainCell $tg 10 6 "?"
Set-PlainCell $tg 10 7 $wildcardChars
Set-PlainCell $tg 10 8 $null
Set-PlainCell $tg 10 9 "a-z,A-Z,0-9"
Set-PlainCell $tg 10 10 $null

# Row 11 - Mode
Set-PlainCell $tg 11 2 "Mode"
Set-PlainCell $tg 11 3 "RobotoCondensed-Regular.ttf"
Set-PlainCell $tg 11 4 24
Set-PlainCell $tg 11 5 4
Set-PlainCell $tg 11 6 "?"
Set-PlainCell $tg 11 7 $null
Set-PlainCell $tg 11 8 $null
Set-PlainCell $tg 11 9 $null
Set-PlainCell $tg 11 10 $null

# Row 12 - button_labels
Set-PlainCell $tg 12 2 "button_labels"
Set-PlainCell $tg 12 3 "verdana.ttf"
Set-PlainCell $tg 12 4 10
Set-PlainCell $tg 12 5 4
Set-PlainCell $tg 12 6 "?"
Set-PlainCell $tg 12 7 $null
Set-PlainCell $tg 12 8 $null
Set-PlainCell $tg 12 9 $null
Set-PlainCell $tg 12 10 $null

# ----------------------- Translation sheet -----------------------------
$tr = $wb.Worksheets.Item("Translation")

# Row 11 - EnteredText (Display)
Set-PlainCell $tr 11 2 "EnteredText"
Set-PlainCell $tr 11 3 "Display"
Set-PlainCell $tr 11 4 "Left"
Set-PlainCell $tr 11 5 $placeholder
Set-PlainCell $tr 11 6 "LTR"

# Row 12 - NumMode (Mode)
Set-PlainCell $tr 12 2 "NumMode"
Set-PlainCell $tr 12 3 "Mode"
Set-PlainCell $tr 12 4 "Center"
Set-PlainCell $tr 12 5 "ABC"
Set-PlainCell $tr 12 6 "LTR"

# Row 13 - AlphaMode (Mode)
Set-PlainCell $tr 13 2 "AlphaMode"
Set-PlainCell $tr 13 3 "Mode"
Set-PlainCell $tr 13 4 "Center"
Set-PlainTextCell $tr 13 5 "123"
Set-PlainCell $tr 13 6 "LTR"

# Row 14 - SingleUseId12 (Default)
Set-PlainCell $tr 14 2 "SingleUseId12"
Set-PlainCell $tr 14 3 "Default"
Set-PlainCell $tr 14 4 "Left"
Set-PlainCell $tr 14 5 "Nazwa miasta:"
Set-PlainCell $tr 14 6 "LTR"

# Row 15 - SingleUseId11 (button_labels)
Set-PlainCell $tr 15 2 "SingleUseId11"
Set-PlainCell $tr 15 3 "button_labels"
Set-PlainCell $tr 15 4 "Center"
Set-PlainCell $tr 15 5 "OK"
Set-PlainCell $tr 15 6 "LTR"

# Row 16 - SingleUseId10 (button_labels)
Set-PlainCell $tr 16 2 "SingleUseId10"
Set-PlainCell $tr 16 3 "button_labels"
Set-PlainCell $tr 16 4 "Center"
Set-PlainCell $tr 16 5 "Cancel"
Set-PlainCell $tr 16 6 "LTR"

# Row 17 - SingleUseId13 (owm_style_20)
Set-PlainCell $tr 17 2 "SingleUseId13"
Set-PlainCell $tr 17 3 "owm_style_20"
Set-PlainCell $tr 17 4 "Left"
Set-PlainCell $tr 17 5 "wpisz co"
Set-PlainCell $tr 17 6 "LTR"
